$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -7
$ws.Range("F4").Value = 2
$ws.Range("F7").Value = -8
$ws.Range("F10").Value = -3
$ws.Range("F12").Value = -1
